$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace every occurrence of "V. cholerae" with "S. sonnei" (cells B10 and B12)
$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    for ($c = 1; $c -le $used.Columns.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq "V. cholerae") {
            $cell.Value = "S. sonnei"
        }
    }
}

# Update the active window's view: zoom to 174% and move the selection to B13
$ws.Range("B13").Select()
$excel.ActiveWindow.Zoom = 174
